$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the _GoBack bookmark that currently sits at the end of the
#    "Input password." paragraph.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2) "Student select "Course list" menu." -> "Student select "Enroll" button."
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "Student select " + [char]0x201C + "Course list" + [char]0x201D + " menu.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Student select " + [char]0x201C + "Enroll" + [char]0x201D + " button.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Delete the two now-redundant steps that used to follow:
#      "Student select "Enroll" button."   (duplicate of the updated step 2)
#      "Go to the student information page."
#    Both are simple single-run ListParagraph items, so after removing the
#    first one the second slides into the same paragraph index.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute(
    "Student select " + [char]0x201C + "Enroll" + [char]0x201D + " button.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dupParagraph = $r2.Paragraphs(1)
$dupParagraph.Range.Delete()

$r3 = $d.Content
$r3.Find.Execute("Go to the student information page.", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$infoParagraph = $r3.Paragraphs(1)
$infoParagraph.Range.Delete()

# ---------------------------------------------------------------------------
# 4) Rewrite the last step's trailing sentence: swap
#      " list that show on the screen."
#    for
#      " list that show on the "Student Information" page."
#    and move the (re-created) _GoBack bookmark so it sits right before the
#    final period, matching the target markup exactly.
# ---------------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("Check the course which was added in the enrolled course list that show on the screen.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$checkParagraph = $target.Paragraphs(1)

$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
    '<w:pPr>' +
      '<w:pStyle w:val="ListParagraph"/>' +
      '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr>' +
      '<w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
      '<w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>' +
      '<w:t xml:space="preserve">Check the course which was added in the </w:t>' +
    '</w:r>' +
    '<w:r>' +
      '<w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>' +
      '<w:t>enrolled course</w:t>' +
    '</w:r>' +
    '<w:r>' +
      '<w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>' +
      '<w:t xml:space="preserve"> list that show on the </w:t>' +
    '</w:r>' +
    '<w:r>' +
      '<w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>' +
      '<w:t>' + [char]0x201C + 'Student Information' + [char]0x201D + ' page</w:t>' +
    '</w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r>' +
      '<w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>' +
      '<w:t>.</w:t>' +
    '</w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$checkParagraph.Range.InsertXML($newParaXml) | Out-Null
